$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "2021" column (R) mirroring the formatting of the preceding
# "2020" column (Q) for each row, then fill in the new values.

# Row 2 (bottom-border separator row, empty value cell)
$ws.Range("Q2").Copy()
$ws.Range("R2").PasteSpecial(-4122)

# Row 3 (year header)
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021

# Row 4 (population count)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 202551

# Row 5 (percentage)
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 2.9794303052841493

$excel.CutCopyMode = $false

# Move / update the active selection to the newly added cell
$ws.Range("R2").Select()
